$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = "Sidkidat"
$ws.Range("C5").Value = "Valiant"
$ws.Range("C7").Value = "UGV"

$ws.Range("D17").Select()
